{"js": "// Word JS API (Office.js) equivalent of the OOXML diff:\n//  - The first paragraph's pPr gains a <w:pBdr> (top/left/bottom/right,\n//    each carrying only w:space=\"5\") and its left indent moves from\n//    120 to 225 twips (6pt -> 11.25pt).\n//  - The placeholder-id run text changes from\n//    \"**ID__AFFARS_5316_topic_8__ID**\" to\n//    \"**ID__AFFARS_SUBPART_5316_206__ID**\", and the trailing run that\n//    held just a single space is removed (the id run ends up as the\n//    only run left in the paragraph).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\n\n// Word's object model always stamps an explicit border line style\n// (w:val) onto a pBdr edge the moment any of its Borders properties are\n// touched, so a pBdr made up of only w:space (no w:val/w:sz/w:color)\n// can't be produced via paragraph.borders. Rebuilding the paragraph\n// from a literal OOXML fragment is the only way to get that exact\n// shape, and it lets the new indent / replacement text be set in the\n// same shot.\nconst targetRange = firstParagraph.getRange();\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"utf-8\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:pPr>' +\n  '<w:pBdr>' +\n  '<w:top w:space=\"5\"/>' +\n  '<w:left w:space=\"5\"/>' +\n  '<w:bottom w:space=\"5\"/>' +\n  '<w:right w:space=\"5\"/>' +\n  '</w:pBdr>' +\n  '<w:spacing w:after=\"0\"/>' +\n  '<w:ind w:left=\"225\"/>' +\n  '<w:jc w:val=\"left\"/>' +\n  '</w:pPr>' +\n  '<w:r>' +\n  '<w:rPr>' +\n  '<w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\"/>' +\n  '<w:b w:val=\"false\"/>' +\n  '<w:i w:val=\"false\"/>' +\n  '<w:color w:val=\"000000\"/>' +\n  '<w:sz w:val=\"22\"/>' +\n  '</w:rPr>' +\n  '<w:t>**ID__AFFARS_SUBPART_5316_206__ID**</w:t>' +\n  '</w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\ntargetRange.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$p = $d.Paragraphs(1)\n\n# Add a paragraph border (top/left/bottom/right) with 5 twip spacing,\n# and update the left indent from 120 -> 225 twips (11.25 pt).\n$pf = $p.Range.ParagraphFormat\n$pf.Borders.DistanceFromTop = 5\n$pf.Borders.DistanceFromLeft = 5\n$pf.Borders.DistanceFromBottom = 5\n$pf.Borders.DistanceFromRight = 5\n$pf.LeftIndent = 11.25\n\n# Replace the paragraph's text (the placeholder id run + trailing space\n# run) with a single updated placeholder id, leaving the paragraph mark\n# - and therefore the paragraph itself - intact.\n$r = $p.Range\n$r.MoveEnd(1, -1) | Out-Null\n$r.Text = \"**ID__AFFARS_SUBPART_5316_206__ID**\"\n"}
